$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the old test rows 8-16 (final data only spans rows 1-7) ---
$ws.Range("A8:B16").EntireRow.Delete()

# --- Remove the hyperlink that used to sit on A7 ---
$ws.Hyperlinks.Delete()
$ws.Range("A7").ClearFormats()

# --- Rewrite the sample login/password rows ---
# Row 2 (987654321 / admin) stays the same value-wise.
# Introduce "admin12" before any of the new " " values so the shared-string
# table is built in the same order as the target workbook.
$ws.Range("B4").Value = "admin12"
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("B5").Value = " "
$ws.Range("A4").Value = 9876543211
$ws.Range("A5").Value = 9876543211
$ws.Range("A6").Value = ""
$ws.Range("A7").Value = 99865845367

# --- Re-apply the left-aligned number style to column A data rows ---
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A7").HorizontalAlignment = -4131

# --- The Hyperlink named cell style is no longer used anywhere; drop it ---
$wb.Styles("Hyperlink").Delete()

# --- Match the final selection recorded in the workbook ---
$ws.Range("C5").Select()
